{"js": "const replacements = [\n  [\"2025-11-11 Tuesday\", \"2025-11-12 Wednesday\"],\n  [\"537\u00d76=3222\", \"329\u00d77=2303\"],\n  [\"459\u00d72=918\", \"725\u00d75=3625\"],\n  [\"629\u00d76=3774\", \"678\u00d73=2034\"],\n  [\"593\u00d72=1186\", \"310\u00d78=2480\"],\n  [\"218\u00d72=436\", \"419\u00d72=838\"],\n  [\"640\u00d73=1920\", \"661\u00d72=1322\"],\n  [\"510\u00d76=3060\", \"301\u00d74=1204\"],\n  [\"972\u00d75=4860\", \"695\u00d75=3475\"],\n  [\"801\u00d72=1602\", \"979\u00d78=7832\"],\n  [\"746\u00d78=5968\", \"745\u00d79=6705\"],\n  [\"291\u00d78=2328\", \"872\u00d73=2616\"],\n  [\"773\u00d76=4638\", \"393\u00d72=786\"],\n  [\"110\u00d75=550\", \"115\u00d75=575\"],\n  [\"772\u00d79=6948\", \"593\u00d74=2372\"],\n  [\"484\u00d77=3388\", \"901\u00d79=8109\"],\n  [\"982\u00d76=5892\", \"431\u00d77=3017\"],\n  [\"505\u00d72=1010\", \"571\u00d74=2284\"],\n  [\"498\u00d72=996\", \"194\u00d78=1552\"],\n  [\"305\u00d74=1220\", \"816\u00d73=2448\"],\n  [\"346\u00d75=1730\", \"596\u00d76=3576\"],\n  [\"752\u00d77=5264\", \"804\u00d72=1608\"],\n  [\"588\u00d78=4704\", \"694\u00d73=2082\"],\n  [\"110\u00d72=220\", \"988\u00d75=4940\"],\n  [\"707\u00d78=5656\", \"987\u00d73=2961\"],\n  [\"275\u00d76=1650\", \"644\u00d75=3220\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  ,@(\"2025-11-11 Tuesday\", \"2025-11-12 Wednesday\")\n  ,@(\"537\u00d76=3222\", \"329\u00d77=2303\")\n  ,@(\"459\u00d72=918\", \"725\u00d75=3625\")\n  ,@(\"629\u00d76=3774\", \"678\u00d73=2034\")\n  ,@(\"593\u00d72=1186\", \"310\u00d78=2480\")\n  ,@(\"218\u00d72=436\", \"419\u00d72=838\")\n  ,@(\"640\u00d73=1920\", \"661\u00d72=1322\")\n  ,@(\"510\u00d76=3060\", \"301\u00d74=1204\")\n  ,@(\"972\u00d75=4860\", \"695\u00d75=3475\")\n  ,@(\"801\u00d72=1602\", \"979\u00d78=7832\")\n  ,@(\"746\u00d78=5968\", \"745\u00d79=6705\")\n  ,@(\"291\u00d78=2328\", \"872\u00d73=2616\")\n  ,@(\"773\u00d76=4638\", \"393\u00d72=786\")\n  ,@(\"110\u00d75=550\", \"115\u00d75=575\")\n  ,@(\"772\u00d79=6948\", \"593\u00d74=2372\")\n  ,@(\"484\u00d77=3388\", \"901\u00d79=8109\")\n  ,@(\"982\u00d76=5892\", \"431\u00d77=3017\")\n  ,@(\"505\u00d72=1010\", \"571\u00d74=2284\")\n  ,@(\"498\u00d72=996\", \"194\u00d78=1552\")\n  ,@(\"305\u00d74=1220\", \"816\u00d73=2448\")\n  ,@(\"346\u00d75=1730\", \"596\u00d76=3576\")\n  ,@(\"752\u00d77=5264\", \"804\u00d72=1608\")\n  ,@(\"588\u00d78=4704\", \"694\u00d73=2082\")\n  ,@(\"110\u00d72=220\", \"988\u00d75=4940\")\n  ,@(\"707\u00d78=5656\", \"987\u00d73=2961\")\n  ,@(\"275\u00d76=1650\", \"644\u00d75=3220\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n  #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n  # Wrap:=1 (wdFindContinue), Replace:=2 (wdReplaceAll)\n  $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}"}
